$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9851366877555847
$ws.Range("B1").Value = 3.372609376907349
$ws.Range("C1").Value = 2.692147016525269
$ws.Range("D1").Value = 0.9088799357414246
$ws.Range("E1").Value = 0.7262665629386902
